$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 4079
$ws.Cells.Item(2, 3).Value = 3
$ws.Cells.Item(2, 4).Value = 430
$ws.Cells.Item(2, 5).Value = 981
$ws.Cells.Item(2, 6).Value = 48
$ws.Cells.Item(2, 7).Value = 3
$ws.Cells.Item(2, 8).Value = 1421
$ws.Cells.Item(2, 9).Value = 21701
$ws.Cells.Item(2, 10).Value = 9
$ws.Cells.Item(2, 11).Value = 497
$ws.Cells.Item(2, 12).Value = 14496

$ws.Cells.Item(3, 2).Value = 821
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(3, 4).Value = 852
$ws.Cells.Item(3, 5).Value = 867
$ws.Cells.Item(3, 6).Value = 5
$ws.Cells.Item(3, 7).Value = 10
$ws.Cells.Item(3, 8).Value = 58
$ws.Cells.Item(3, 9).Value = 7568
$ws.Cells.Item(3, 10).Value = 14
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 0

$ws.Cells.Item(4, 2).Value = 265
$ws.Cells.Item(4, 3).Value = 0
$ws.Cells.Item(4, 4).Value = 133
$ws.Cells.Item(4, 5).Value = 142
$ws.Cells.Item(4, 6).Value = 8
$ws.Cells.Item(4, 7).Value = 1
$ws.Cells.Item(4, 8).Value = 567
$ws.Cells.Item(4, 9).Value = 4846
$ws.Cells.Item(4, 10).Value = 2
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 0

$ws.Cells.Item(5, 2).Value = 803
$ws.Cells.Item(5, 3).Value = 0
$ws.Cells.Item(5, 4).Value = 377
$ws.Cells.Item(5, 5).Value = 377
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 8061
$ws.Cells.Item(5, 10).Value = 11
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 0

$ws.Cells.Item(6, 2).Value = 266
$ws.Cells.Item(6, 3).Value = 0
$ws.Cells.Item(6, 4).Value = 146
$ws.Cells.Item(6, 5).Value = 151
$ws.Cells.Item(6, 6).Value = 4
$ws.Cells.Item(6, 7).Value = 1
$ws.Cells.Item(6, 8).Value = 267
$ws.Cells.Item(6, 9).Value = 5105
$ws.Cells.Item(6, 10).Value = 1
$ws.Cells.Item(6, 11).Value = 0
$ws.Cells.Item(6, 12).Value = 0

$ws.Cells.Item(7, 2).Value = 262
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(7, 4).Value = 3
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 0
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 12).Value = 0

$ws.Cells.Item(8, 2).Value = 1
$ws.Cells.Item(8, 3).Value = 5
$ws.Cells.Item(8, 4).Value = 244
$ws.Cells.Item(8, 5).Value = 306
$ws.Cells.Item(8, 6).Value = 46
$ws.Cells.Item(8, 7).Value = 11
$ws.Cells.Item(8, 8).Value = 1559
$ws.Cells.Item(8, 9).Value = 2397
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(8, 11).Value = 0
$ws.Cells.Item(8, 12).Value = 0

$ws.Cells.Item(9, 2).Value = 777
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(9, 4).Value = 56
$ws.Cells.Item(9, 5).Value = 56
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(9, 7).Value = 0
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 9).Value = 7857
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(9, 12).Value = 0

$ws.Cells.Item(10, 2).Value = 766
$ws.Cells.Item(10, 3).Value = 0
$ws.Cells.Item(10, 4).Value = 58
$ws.Cells.Item(10, 5).Value = 61
$ws.Cells.Item(10, 6).Value = 3
$ws.Cells.Item(10, 7).Value = 0
$ws.Cells.Item(10, 8).Value = 1611
$ws.Cells.Item(10, 9).Value = 17118
$ws.Cells.Item(10, 10).Value = 0
$ws.Cells.Item(10, 11).Value = 0
$ws.Cells.Item(10, 12).Value = 0

$ws.Cells.Item(11, 2).Value = 984
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(11, 4).Value = 213
$ws.Cells.Item(11, 5).Value = 214
$ws.Cells.Item(11, 6).Value = 2
$ws.Cells.Item(11, 7).Value = 0
$ws.Cells.Item(11, 8).Value = 1985
$ws.Cells.Item(11, 9).Value = 57005
$ws.Cells.Item(11, 10).Value = 6
$ws.Cells.Item(11, 11).Value = 0
$ws.Cells.Item(11, 12).Value = 0

$ws.Cells.Item(12, 2).Value = 829
$ws.Cells.Item(12, 3).Value = 0
$ws.Cells.Item(12, 4).Value = 38
$ws.Cells.Item(12, 5).Value = 37
$ws.Cells.Item(12, 6).Value = 0
$ws.Cells.Item(12, 7).Value = 0
$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(12, 9).Value = 30000
$ws.Cells.Item(12, 10).Value = 2
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 12).Value = 0

$ws.Cells.Item(13, 2).Value = 575
$ws.Cells.Item(13, 3).Value = 0
$ws.Cells.Item(13, 4).Value = 246
$ws.Cells.Item(13, 5).Value = 325
$ws.Cells.Item(13, 6).Value = 84
$ws.Cells.Item(13, 7).Value = 1
$ws.Cells.Item(13, 8).Value = 4560
$ws.Cells.Item(13, 9).Value = 4028
$ws.Cells.Item(13, 10).Value = 5
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(13, 12).Value = 0

$ws.Cells.Item(14, 2).Value = 804
$ws.Cells.Item(14, 3).Value = 0
$ws.Cells.Item(14, 4).Value = 269
$ws.Cells.Item(14, 5).Value = 272
$ws.Cells.Item(14, 6).Value = 2
$ws.Cells.Item(14, 7).Value = 1
$ws.Cells.Item(14, 8).Value = 74
$ws.Cells.Item(14, 9).Value = 8539
$ws.Cells.Item(14, 10).Value = 5
$ws.Cells.Item(14, 11).Value = 0
$ws.Cells.Item(14, 12).Value = 0

$ws.Cells.Item(15, 2).Value = 264
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(15, 4).Value = 74
$ws.Cells.Item(15, 5).Value = 77
$ws.Cells.Item(15, 6).Value = 3
$ws.Cells.Item(15, 7).Value = 0
$ws.Cells.Item(15, 8).Value = 390
$ws.Cells.Item(15, 9).Value = 5000
$ws.Cells.Item(15, 10).Value = 1
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(15, 12).Value = 0

$ws.Cells.Item(16, 2).Value = 92
$ws.Cells.Item(16, 3).Value = 2
$ws.Cells.Item(16, 4).Value = 424
$ws.Cells.Item(16, 5).Value = 615
$ws.Cells.Item(16, 6).Value = 38
$ws.Cells.Item(16, 7).Value = 4
$ws.Cells.Item(16, 8).Value = 3501
$ws.Cells.Item(16, 9).Value = 10764
$ws.Cells.Item(16, 10).Value = 6
$ws.Cells.Item(16, 11).Value = 147
$ws.Cells.Item(16, 12).Value = 3349

$ws.Cells.Item(17, 2).Value = 62
$ws.Cells.Item(17, 3).Value = 0
$ws.Cells.Item(17, 4).Value = 26
$ws.Cells.Item(17, 5).Value = 27
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(17, 7).Value = 0
$ws.Cells.Item(17, 8).Value = 0
$ws.Cells.Item(17, 9).Value = 5385
$ws.Cells.Item(17, 10).Value = 0
$ws.Cells.Item(17, 11).Value = 1
$ws.Cells.Item(17, 12).Value = 370

$ws.Cells.Item(18, 2).Value = 223
$ws.Cells.Item(18, 3).Value = 0
$ws.Cells.Item(18, 4).Value = 529
$ws.Cells.Item(18, 5).Value = 627
$ws.Cells.Item(18, 6).Value = 25
$ws.Cells.Item(18, 7).Value = 1
$ws.Cells.Item(18, 8).Value = 399
$ws.Cells.Item(18, 9).Value = 5501
$ws.Cells.Item(18, 10).Value = 10
$ws.Cells.Item(18, 11).Value = 73
$ws.Cells.Item(18, 12).Value = 1164

$ws.Cells.Item(19, 2).Value = 778
$ws.Cells.Item(19, 3).Value = 1
$ws.Cells.Item(19, 4).Value = 195
$ws.Cells.Item(19, 5).Value = 203
$ws.Cells.Item(19, 6).Value = 5
$ws.Cells.Item(19, 7).Value = 1
$ws.Cells.Item(19, 8).Value = 284
$ws.Cells.Item(19, 9).Value = 5018
$ws.Cells.Item(19, 10).Value = 1
$ws.Cells.Item(19, 11).Value = 0
$ws.Cells.Item(19, 12).Value = 0

$ws.Cells.Item(20, 2).Value = 1396
$ws.Cells.Item(20, 3).Value = 0
$ws.Cells.Item(20, 4).Value = 60
$ws.Cells.Item(20, 5).Value = 78
$ws.Cells.Item(20, 6).Value = 18
$ws.Cells.Item(20, 7).Value = 0
$ws.Cells.Item(20, 8).Value = 8800
$ws.Cells.Item(20, 9).Value = 14576
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 11).Value = 0
$ws.Cells.Item(20, 12).Value = 0

$ws.Cells.Item(21, 2).Value = 295
$ws.Cells.Item(21, 3).Value = 1
$ws.Cells.Item(21, 4).Value = 63
$ws.Cells.Item(21, 5).Value = 70
$ws.Cells.Item(21, 6).Value = 6
$ws.Cells.Item(21, 7).Value = 0
$ws.Cells.Item(21, 8).Value = 1250
$ws.Cells.Item(21, 9).Value = 13770
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 11).Value = 0
$ws.Cells.Item(21, 12).Value = 0

$ws.Cells.Item(22, 2).Value = 187
$ws.Cells.Item(22, 3).Value = 0
$ws.Cells.Item(22, 4).Value = 1
$ws.Cells.Item(22, 5).Value = 1
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(22, 7).Value = 0
$ws.Cells.Item(22, 8).Value = 0
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 0
$ws.Cells.Item(22, 12).Value = 0

$ws.Cells.Item(23, 2).Value = 813
$ws.Cells.Item(23, 3).Value = 0
$ws.Cells.Item(23, 4).Value = 15
$ws.Cells.Item(23, 5).Value = 16
$ws.Cells.Item(23, 6).Value = 1
$ws.Cells.Item(23, 7).Value = 0
$ws.Cells.Item(23, 8).Value = 625
$ws.Cells.Item(23, 9).Value = 8000
$ws.Cells.Item(23, 10).Value = 0
$ws.Cells.Item(23, 11).Value = 0
$ws.Cells.Item(23, 12).Value = 0

$ws.Cells.Item(24, 2).Value = 988
$ws.Cells.Item(24, 3).Value = 0
$ws.Cells.Item(24, 4).Value = 43
$ws.Cells.Item(24, 5).Value = 48
$ws.Cells.Item(24, 6).Value = 4
$ws.Cells.Item(24, 7).Value = 1
$ws.Cells.Item(24, 8).Value = 1710
$ws.Cells.Item(24, 9).Value = 16002
$ws.Cells.Item(24, 10).Value = 0
$ws.Cells.Item(24, 11).Value = 0
$ws.Cells.Item(24, 12).Value = 0

$ws.Cells.Item(25, 2).Value = 74
$ws.Cells.Item(25, 3).Value = 0
$ws.Cells.Item(25, 4).Value = 64
$ws.Cells.Item(25, 5).Value = 85
$ws.Cells.Item(25, 6).Value = 4
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(25, 8).Value = 471
$ws.Cells.Item(25, 9).Value = 4063
$ws.Cells.Item(25, 10).Value = 5
$ws.Cells.Item(25, 11).Value = 17
$ws.Cells.Item(25, 12).Value = 20

$ws.Cells.Item(26, 2).Value = 133
$ws.Cells.Item(26, 3).Value = 0
$ws.Cells.Item(26, 4).Value = 0
$ws.Cells.Item(26, 5).Value = 0
$ws.Cells.Item(26, 6).Value = 0
$ws.Cells.Item(26, 7).Value = 0
$ws.Cells.Item(26, 8).Value = 0
$ws.Cells.Item(26, 9).Value = 0
$ws.Cells.Item(26, 10).Value = 0
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 12).Value = 0

$ws.Cells.Item(27, 2).Value = 6
$ws.Cells.Item(27, 3).Value = 2
$ws.Cells.Item(27, 4).Value = 106
$ws.Cells.Item(27, 5).Value = 118
$ws.Cells.Item(27, 6).Value = 9
$ws.Cells.Item(27, 7).Value = 1
$ws.Cells.Item(27, 8).Value = 769
$ws.Cells.Item(27, 9).Value = 1810
$ws.Cells.Item(27, 10).Value = 4
$ws.Cells.Item(27, 11).Value = 0
$ws.Cells.Item(27, 12).Value = 0

